{"js": "// Map of old text -> new text for this edit. Each old string is unique\n// within the document, so a simple exact search/replace per pair is safe.\nconst replacements = [\n  [\"2025-04-24 Thursday\", \"2025-04-25 Friday\"],\n  [\"186\u00f76=\", \"291\u00f79=\"],\n  [\"326\u00f79=\", \"897\u00f76=\"],\n  [\"960\u00f74=\", \"935\u00f75=\"],\n  [\"405\u00f74=\", \"122\u00f73=\"],\n  [\"627\u00f79=\", \"643\u00f78=\"],\n  [\"708\u00f76=\", \"324\u00f78=\"],\n  [\"474\u00f76=\", \"208\u00f76=\"],\n  [\"620\u00f78=\", \"841\u00f73=\"],\n  [\"489\u00f74=\", \"646\u00f72=\"],\n  [\"521\u00f79=\", \"475\u00f76=\"],\n  [\"304\u00f79=\", \"162\u00f77=\"],\n  [\"902\u00f73=\", \"874\u00f79=\"],\n  [\"251\u00f73=\", \"472\u00f78=\"],\n  [\"471\u00f74=\", \"140\u00f72=\"],\n  [\"845\u00f79=\", \"729\u00f74=\"],\n  [\"802\u00f78=\", \"157\u00f73=\"],\n  [\"515\u00f78=\", \"938\u00f72=\"],\n  [\"939\u00f72=\", \"611\u00f73=\"],\n  [\"201\u00f76=\", \"906\u00f73=\"],\n  [\"728\u00f77=\", \"696\u00f75=\"],\n  [\"298\u00f76=\", \"130\u00f78=\"],\n  [\"278\u00f78=\", \"695\u00f74=\"],\n  [\"576\u00f74=\", \"847\u00f75=\"],\n  [\"583\u00f79=\", \"826\u00f78=\"],\n  [\"404\u00f72=\", \"358\u00f78=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Map of old text -> new text for this edit. Each old string is unique\n# within the document, so a simple Find/Replace pair per entry is safe.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-04-24 Thursday\", \"2025-04-25 Friday\"),\n    @(\"186\u00f76=\", \"291\u00f79=\"),\n    @(\"326\u00f79=\", \"897\u00f76=\"),\n    @(\"960\u00f74=\", \"935\u00f75=\"),\n    @(\"405\u00f74=\", \"122\u00f73=\"),\n    @(\"627\u00f79=\", \"643\u00f78=\"),\n    @(\"708\u00f76=\", \"324\u00f78=\"),\n    @(\"474\u00f76=\", \"208\u00f76=\"),\n    @(\"620\u00f78=\", \"841\u00f73=\"),\n    @(\"489\u00f74=\", \"646\u00f72=\"),\n    @(\"521\u00f79=\", \"475\u00f76=\"),\n    @(\"304\u00f79=\", \"162\u00f77=\"),\n    @(\"902\u00f73=\", \"874\u00f79=\"),\n    @(\"251\u00f73=\", \"472\u00f78=\"),\n    @(\"471\u00f74=\", \"140\u00f72=\"),\n    @(\"845\u00f79=\", \"729\u00f74=\"),\n    @(\"802\u00f78=\", \"157\u00f73=\"),\n    @(\"515\u00f78=\", \"938\u00f72=\"),\n    @(\"939\u00f72=\", \"611\u00f73=\"),\n    @(\"201\u00f76=\", \"906\u00f73=\"),\n    @(\"728\u00f77=\", \"696\u00f75=\"),\n    @(\"298\u00f76=\", \"130\u00f78=\"),\n    @(\"278\u00f78=\", \"695\u00f74=\"),\n    @(\"576\u00f74=\", \"847\u00f75=\"),\n    @(\"583\u00f79=\", \"826\u00f78=\"),\n    @(\"404\u00f72=\", \"358\u00f78=\")\n)\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Forward = $true\n    $find.Wrap = 1  # wdFindContinue\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute([ref]$old, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]1, [ref]$false, [ref]$new, [ref]2) | Out-Null\n}\n"}
